$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, shifting existing rows 59-64 down to 60-65
$ws.Rows.Item(59).Insert()

# Copy style (number format) of the date cell in column D from the row below (now row 60) to the new row 59
$ws.Range("D60").Copy()
$ws.Range("D59").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate new row 59 with data
$ws.Range("A59").Value = 3
$ws.Range("B59").Value = "Femacal de La Calera"
$ws.Range("C59").Value = "Coquimbo"
$ws.Range("D59").Value = 44918
$ws.Range("E59").Value = 5
$ws.Range("F59").Value = 300000000
$ws.Range("G59").Value = "Espárragos"
$ws.Range("H59").Value = "Verde"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 550
$ws.Range("K59").Value = 1500
$ws.Range("L59").Value = 1500
$ws.Range("M59").Value = 1500
$ws.Range("N59").Value = "$/kilo"
$ws.Range("O59").Value = "Provincia de Quillota"
$ws.Range("P59").Value = 1500
$ws.Range("Q59").Value = 1
$ws.Range("R59").Value = "Hortaliza"
